$d = $word.ActiveDocument

function Add-Para {
    $count = $d.Paragraphs.Count
    $p = $d.Paragraphs.Item($count)
    $r = $p.Range
    $r.InsertParagraphAfter()
    $newCount = $d.Paragraphs.Count
    $newP = $d.Paragraphs.Item($newCount)
    return $newP
}

function Add-Run {
    param($para, [string]$text, [bool]$italic = $false, [bool]$bold = $false)
    $preEnd = $para.Range.End - 1
    $ir = $d.Range($preEnd, $preEnd)
    $ir.InsertAfter($text)
    $postEnd = $para.Range.End - 1
    $runRange = $d.Range($preEnd, $postEnd)
    if ($italic) { $runRange.Italic = 1 }
    if ($bold) { $runRange.Bold = 1 }
}

# Insert the new "Knärot" section (heading, body paragraphs, and references)
# right after the "BILAGA 1 - Fridlysta arter" heading, before the section break.
# --- Paragraph 0 ---
$p = Add-Para
$p.Style = "Heading 1"
Add-Run $p "Knärot – ekologi samt krav på livsmiljön" $false $false

# --- Paragraph 1 ---
$p = Add-Para
$p.Style = "Normal"
Add-Run $p "Knärot är fridlyst enligt 8 och 15 §§ artskyddsförordningen och klassad som sårbar (VU) enligt rödlistan 2020. Knärot är beroende av hög och jämn luftfuktighet i gamla, ostörda skogsmiljöer och är känslig för snabba förändringar av ljus-/vindförhållanden eller uttorkning. På grund av ett alltför intensivt skogsbruk har den minskat med 40 (25-50) % under de senaste 60 åren och i framtiden bedöms minskningstakten uppgå till 30 (20-40) %. Till följd av att arten har en dokumenterat högre minskningstakt iförhållande till sin generationstid än vad som tidigare varit känt (data från Riksskogstaxeringen) höjdes den till hotkategori sårbar (VU) i rödlistan 2020 (Artdatabanken, 2021)." $false $false

# --- Paragraph 2 ---
$p = Add-Para
$p.Style = "Normal"
Add-Run $p "Samuel Johnsons doktorsavhandling " $false $false
Add-Run $p "“Retention Forestry as a Conservation Measure for Boreal Forest Ground Vegetation“" $true $false
Add-Run $p " (SLU, Uppsala 2014) visar att det krävs väl tilltagna skyddszoner för att knärotens växtplatser inte ska ta skada av skogsbruksåtgärder i intilliggande områden: " $false $false
Add-Run $p "“Study III shows that retention patches smaller than 0.5 ha do not lifeboat the sensitive forest herb G. repens, a species that depend on stable microclimatic conditions typical for intact forest stands.” " $true $false
Add-Run $p "Vidare " $false $false
Add-Run $p "“More sensitive forest species are not lifeboated in retention patches ranging from 0.05 to 0.5 ha (Papers II & III).”" $true $false

# --- Paragraph 3 ---
$p = Add-Para
$p.Style = "Normal"
Add-Run $p "Johnsons (2014) rekommendation på minst 50 meters breda skyddszoner runt knärotens växtplatser motsvarar en areal på 0,78 hektar, vilket ligger i linje med andra studier som gjorts på känsliga skogsarter: " $false $false
Add-Run $p "“In study III I also show that translocated specimens of G. repens survives well in mature forests at least 50 m from the nearest edge to an open area. Moreover, measures of temperature and humidity show that such distances from an open area is far enough to offer a microclimate that is more stable compared to what present in retention patches of around 0.1 ha. This means that the very centre of a circular patch with radius 50 m (equals a size of 0.78 ha) should offer conditions similar to interior forest and would perhaps be a suitable habitat for G. repens and similar species. Previous studies from both North America and Sweden have also concluded that patches between 0.5 and one ha are sufficient for preserving interior forest vegetation as well as sensitive lichens and bryophytes (de Graaf & Roberts 2009; Halpern et al. 2012; Rudolphi et al. 2014).”" $true $false

# --- Paragraph 4 ---
$p = Add-Para
$p.Style = "Normal"
Add-Run $p "En nyligen publicerad vetenskaplig uppsats av Koelmeijer m.fl. (2022) inkluderar orkidén knärots skyddsbehov. I uppsatsen berörs problemet med uttorkning för växter, bl.a. för knärot, ett problem som blivit accentuerat på grund av den pågående klimatförändringen och torra somrar, t.ex. den exceptionellt torra sommaren 2018. I uppsatsen undersöks områden med tre olika avstånd från kalhyggeskant med avseende på skydd bl.a. för knärot. Det första området har avstånd upp till 20 m från hyggeskant (Strong edge effect), det andra 20 – 40 m från hyggeskant (Weak edge effect) och det tredje avser större avstånd från hyggeskant, där kanteffekten anses vara försumbar (Interior). Ett resultat var att man fann stor eller mycket stor uttorkningseffekt på känsliga och rödlistade skogsarter vid de kortare avstånden till hyggeskant, medan effekt av uttorkning inte konstaterades på större avstånd (Interior). För orkidén knärot fann man en rik förekomst (upp till 0,06 dm2/m2) på stort avstånd från hyggeskant (Interior), medan förekomsten var liten eller närmast försumbar i de områden som klassificerades som Weak edge effect respektive Strong edge effect. Arbetet påpekar att de allt oftare förekommande torra somrarna ger ytterligare skäl att utöka skyddsavståndet från hyggen till den fuktkrävande arten knärot (Koelmeijer m.fl., 2022)." $false $false

# --- Paragraph 5 ---
$p = Add-Para
$p.Style = "Normal"
Add-Run $p "Även Skogsstyrelsens egen vägledning för hänsyn till knärot ligger i linje med ovanstående forskningsstudier. Av vägledningen framgår det att för med hög sannolikhet kunna bevara befintliga förekomster krävs relativt stora avsättningar av uppvuxen skog med slutet och relativt tätt kronskikt. Som riktlinje kan krävas ett avstånd på 50 meter in från brynet för att vidmakthålla ett fungerande mikroklimat. Detta innebär att fristående hänsynsytor för många arter (kärlväxter, lavar och mossor) kan behöva ha en area överstigande 0,8 hektar (cirkelyta med radien 50 meter = 0,78 hektar) för att bibehålla lokalklimatet. Även ganska små förändringar i form av förändrade ljus- och fuktighetsförhållanden, till exempel till följd av gallring, kan leda till att arten försvinner till följd av konkurrens med mera ljuskrävande och snabbväxande arter (Skogsstyrelsen, 2022)." $false $false

# --- Paragraph 6 ---
$p = Add-Para
$p.Style = "Heading 2"
Add-Run $p "Referenser - knärot" $false $false

# --- Paragraph 7 ---
$p = Add-Para
$p.Style = "Normal"
Add-Run $p "de Graaf M & Roberts M.R., 2009. " $false $false
Add-Run $p "Short-term response of the herbaceous layer within leave patches after harvest. " $true $false
Add-Run $p "Forest Ecology and Management 257, 1014-1025" $false $false

# --- Paragraph 8 ---
$p = Add-Para
$p.Style = "Normal"
Add-Run $p "Halpern, C. B., Halaj, J., Evans, S. A., & Dovciak, M., 2012. " $false $false
Add-Run $p "Level and pattern of overstory retention interact to shape long-term responses of understories to timber harvest. " $true $false
Add-Run $p "Ecological Applications, 22, 2049-2064 " $false $false

# --- Paragraph 9 ---
$p = Add-Para
$p.Style = "Normal"
Add-Run $p "Koelmeijer, I. A., Ehrlén, J., Jönsson, M., De Frenne, P., Berg, P., Andersson, J., Weibull, H. & Hylander, N. 2022. " $false $false
Add-Run $p "Interactive effects of drought and edge exposure on old-growth forest understory species. " $true $false
Add-Run $p "Landscape Ecology, 37, sid 1839-1853" $false $false

# --- Paragraph 10 ---
$p = Add-Para
$p.Style = "Normal"
Add-Run $p "Rudolphi, J., Jönsson, M. T., & Gustafsson, L., 2014. " $false $false
Add-Run $p "Biological legacies buffer local species extinction after logging. " $true $false
Add-Run $p "Journal of Applied Ecology. 51, 53-62." $false $false

# --- Paragraph 11 ---
$p = Add-Para
$p.Style = "Normal"
Add-Run $p "Skogsstyrelsen, 2022. " $false $false
Add-Run $p "Vägledning för hänsyn till knärot. " $true $false
Add-Run $p "https://www.skogsstyrelsen.se/lag-och-tillsyn/artskydd/vagledningar-och-kunskapsstod-artskydd/vagledning-for-hansyn-till-knarot/" $false $false

# --- Paragraph 12 ---
$p = Add-Para
$p.Style = "Normal"
Add-Run $p "SLU Artdatabanken, 2021. " $false $false
Add-Run $p "Artfaktablad. Naturvård – artfakta. " $true $false
Add-Run $p "SLU Artdatabanken, Uppsala " $false $false

# Update the date in the first-page header from 2023-09-13 to 2023-09-15.
$sec = $d.Sections.Item(1)
$hdr = $sec.Headers.Item(2)
$hdrRange = $hdr.Range
$hdrRange.Find.Execute("2023-09-13", $true, $false, $false, $false, $false, $true, 1, $false, "2023-09-15", 2)
